# semana 33 de 2025
# Adds the "week 33" column (AJ) to the weekly IRA/UCI revision sheet:
#  - AJ1 gets the new header label "33" (text, matching the other week
#    headers in row 1 which are stored as text, not numbers)
#  - AJ2..AJ58 get the week-33 counts for every UPGD row that already had
#    a full run of weekly data (rows that are missing data for other
#    weeks are left untouched, matching the source diff)
#  - AH28 (week 31 for "CLINICA LOS ROSALES") is corrected from 103 to 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AJ is column 36
$col = 36

# Header cell: must be stored as TEXT "33", like D1:AI1 ("1".."32"),
# not as a number. A leading apostrophe forces Excel to keep it as text
# while still only displaying "33".
$ws.Cells.Item(1, $col).Value = "'33"

# Week-33 values for each data row (row number -> value). Rows not
# listed here have no week-33 figure reported and are left as-is.
$weekValues = [ordered]@{
    2  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    13 = 0
    14 = 0
    15 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    29 = 1
    31 = 0
    34 = 0
    35 = 3
    36 = 0
    37 = 0
    38 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    58 = 0
}

foreach ($row in $weekValues.Keys) {
    $ws.Cells.Item($row, $col).Value = $weekValues[$row]
}

# Correction to an existing value: AH28 (week 31) was 103, now 0.
$ws.Cells.Item(28, 34).Value = 0
